$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vcan"
$ws.Range("C2").Value = "Egfr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2.0
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.529781
$ws.Range("H2").Value = 4.589343
$ws.Range("I2").Value = 0.01315047351877542
$ws.Range("J2").Value = 0.01315047351877542
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 1.307106666666667
$ws.Range("N2").Value = 3.92132
$ws.Range("O2").Value = 0.01256263154946851
$ws.Range("P2").Value = 0.01256263154946851
$ws.Range("Q2").Value = 1.99958694364
$ws.Range("R2").Value = 17.99628249276
$ws.Range("S2").Value = 0.0001652045535174183
$ws.Range("T2").Value = 0.0001652045535174183

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vcan"
$ws.Range("C3").Value = "Egfr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2.0
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.529781
$ws.Range("H3").Value = 4.589343
$ws.Range("I3").Value = 0.01315047351877542
$ws.Range("J3").Value = 0.01315047351877542
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 80.22623699999998
$ws.Range("N3").Value = 240.678711
$ws.Range("O3").Value = 0.77105616682495
$ws.Range("P3").Value = 0.77105616682495
$ws.Range("Q3").Value = 122.728573064097
$ws.Range("R3").Value = 1104.557157576873
$ws.Range("S3").Value = 0.01013975370331999
$ws.Range("T3").Value = 0.01013975370331999

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vcan"
$ws.Range("C4").Value = "Egfr"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2.0
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.529781
$ws.Range("H4").Value = 4.589343
$ws.Range("I4").Value = 0.01315047351877542
$ws.Range("J4").Value = 0.01315047351877542
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 22.51385866666667
$ws.Range("N4").Value = 67.541576
$ws.Range("O4").Value = 0.2163812016255815
$ws.Range("P4").Value = 0.2163812016255815
$ws.Range("Q4").Value = 34.441273224952
$ws.Range("R4").Value = 309.971459024568
$ws.Range("S4").Value = 0.002845515261938014
$ws.Range("T4").Value = 0.002845515261938014

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Vcan"
$ws.Range("C5").Value = "Egfr"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 103.676216
$ws.Range("H5").Value = 311.028648
$ws.Range("I5").Value = 0.8912330150752564
$ws.Range("J5").Value = 0.8912330150752563
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 1.307106666666667
$ws.Range("N5").Value = 3.92132
$ws.Range("O5").Value = 0.01256263154946851
$ws.Range("P5").Value = 0.01256263154946851
$ws.Range("Q5").Value = 135.5158731083733
$ws.Range("R5").Value = 1219.64285797536
$ws.Range("S5").Value = 0.01119623199311236
$ws.Range("T5").Value = 0.01119623199311236

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Vcan"
$ws.Range("C6").Value = "Egfr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 103.676216
$ws.Range("H6").Value = 311.028648
$ws.Range("I6").Value = 0.8912330150752564
$ws.Range("J6").Value = 0.8912330150752563
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 80.22623699999998
$ws.Range("N6").Value = 240.678711
$ws.Range("O6").Value = 0.77105616682495
$ws.Range("P6").Value = 0.77105616682495
$ws.Range("Q6").Value = 8317.55267607919
$ws.Range("R6").Value = 74857.97408471271
$ws.Range("S6").Value = 0.6871907123517701
$ws.Range("T6").Value = 0.6871907123517701

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Vcan"
$ws.Range("C7").Value = "Egfr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 103.676216
$ws.Range("H7").Value = 311.028648
$ws.Range("I7").Value = 0.8912330150752564
$ws.Range("J7").Value = 0.8912330150752563
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 22.51385866666667
$ws.Range("N7").Value = 67.541576
$ws.Range("O7").Value = 0.2163812016255815
$ws.Range("P7").Value = 0.2163812016255815
$ws.Range("Q7").Value = 2334.151674118805
$ws.Range("R7").Value = 21007.36506706925
$ws.Range("S7").Value = 0.192846070730374
$ws.Range("T7").Value = 0.1928460707303739

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Vcan"
$ws.Range("C8").Value = "Egfr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 11.12297
$ws.Range("H8").Value = 33.36891
$ws.Range("I8").Value = 0.09561651140596822
$ws.Range("J8").Value = 0.09561651140596822
$ws.Range("K8").Value = 3.0
$ws.Range("L8").Value = 1.0
$ws.Range("M8").Value = 1.307106666666667
$ws.Range("N8").Value = 3.92132
$ws.Range("O8").Value = 0.01256263154946851
$ws.Range("P8").Value = 0.01256263154946851
$ws.Range("Q8").Value = 14.53890824013333
$ws.Range("R8").Value = 130.8501741612
$ws.Range("S8").Value = 0.001201195002838732
$ws.Range("T8").Value = 0.001201195002838732

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Vcan"
$ws.Range("C9").Value = "Egfr"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 11.12297
$ws.Range("H9").Value = 33.36891
$ws.Range("I9").Value = 0.09561651140596822
$ws.Range("J9").Value = 0.09561651140596822
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 80.22623699999998
$ws.Range("N9").Value = 240.678711
$ws.Range("O9").Value = 0.77105616682495
$ws.Range("P9").Value = 0.77105616682495
$ws.Range("Q9").Value = 892.3540273638898
$ws.Range("R9").Value = 8031.186246275009
$ws.Range("S9").Value = 0.07372570076985997
$ws.Range("T9").Value = 0.07372570076985997

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Vcan"
$ws.Range("C10").Value = "Egfr"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 11.12297
$ws.Range("H10").Value = 33.36891
$ws.Range("I10").Value = 0.09561651140596822
$ws.Range("J10").Value = 0.09561651140596822
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 22.51385866666667
$ws.Range("N10").Value = 67.541576
$ws.Range("O10").Value = 0.2163812016255815
$ws.Range("P10").Value = 0.2163812016255815
$ws.Range("Q10").Value = 250.4209745335733
$ws.Range("R10").Value = 2253.78877080216
$ws.Range("S10").Value = 0.02068961563326952
$ws.Range("T10").Value = 0.02068961563326952

